$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 45
$ws1.Range("F4").Value = 1550
$ws1.Range("F5").Value = 243
$ws1.Range("F6").Value = 54
$ws1.Range("F7").Value = 1074
$ws1.Range("F8").Value = 10108
$ws1.Range("F10").Value = 130
$ws1.Range("F11").Value = 253
$ws1.Range("F14").Value = 6996
$ws1.Range("F15").Value = 1092
$ws1.Range("F16").Value = 653
$ws1.Range("F18").Value = 222

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 45
$ws4.Range("F4").Value = 1550
$ws4.Range("F5").Value = 243
$ws4.Range("F7").Value = 54
$ws4.Range("F8").Value = 1074
$ws4.Range("F11").Value = 10108
$ws4.Range("F13").Value = 130
$ws4.Range("F14").Value = 253
$ws4.Range("F17").Value = 6996
$ws4.Range("F18").Value = 1092
$ws4.Range("F19").Value = 653
$ws4.Range("F21").Value = 222
